$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Login with valid username and password", "PASSED", "chrome"),
    @("Create a country", "PASSED", "chrome"),
    @("Create and Delete Cities", "PASSED", "chrome"),
    @("Login with valid username and password", "PASSED", "chrome"),
    @("Create a country", "FAILED", "chrome"),
    @("Create and Delete Cities", "FAILED", "chrome"),
    @("Login with valid username and password", "PASSED", "chrome"),
    @("Create a country", "PASSED", "chrome"),
    @("Create and Delete Cities", "PASSED", "chrome")
)

$startRow = 441
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
